$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column AC (29th column) from width 8 to width 7.
# ColumnWidth is expressed in character units; the stored OOXML "width"
# attribute includes Excel's internal padding offset (~0.8333 chars), so we
# back that offset out to land exactly on the target stored width of 7.
$ws.Columns.Item(29).ColumnWidth = 6.166666666666667

# Apply "custom accuracy" - round row 5's measurement columns (B:AH) down
# to two decimal places instead of three.
$ws.Range("B5").Value = 16.81
$ws.Range("C5").Value = 12.28
$ws.Range("D5").Value = 1.13
$ws.Range("E5").Value = 36.26
$ws.Range("F5").Value = 29.97
$ws.Range("G5").Value = 13.23
$ws.Range("H5").Value = 49.96
$ws.Range("I5").Value = 20.36
$ws.Range("J5").Value = 8.98
$ws.Range("K5").Value = 13.43
$ws.Range("L5").Value = 14.66
$ws.Range("M5").Value = 15.36
$ws.Range("N5").Value = 4.23
$ws.Range("O5").Value = 13.16
$ws.Range("P5").Value = 18.68
$ws.Range("Q5").Value = 11.14
$ws.Range("R5").Value = 0.83
$ws.Range("S5").Value = 0.72
$ws.Range("T5").Value = 192.96
$ws.Range("U5").Value = 36.77
$ws.Range("V5").Value = 12.15
$ws.Range("W5").Value = 24.63
$ws.Range("X5").Value = 13.13
$ws.Range("Y5").Value = 1.83
$ws.Range("Z5").Value = 24.33
$ws.Range("AA5").Value = 10.73
$ws.Range("AB5").Value = 9.58
$ws.Range("AC5").Value = 11.23
$ws.Range("AD5").Value = 15.32
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 45.21
$ws.Range("AG5").Value = 6.81
$ws.Range("AH5").Value = 15.19

# Data row 6 (the 1000th-reading addendum) is dropped entirely, which also
# shrinks the sheet dimension from A1:AH6 to A1:AH5.
$ws.Rows.Item(6).Delete()
